$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21, pushing the existing summary rows (old 21-24) down to 22-25.
$ws.Rows("21").Insert()

# Fill in the new data row (2014-02-22, 21:45 -> 22:30) consistent with the other entries.
$ws.Range("A21").Value = 2014
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 0.90625
$ws.Range("E21").Value = 0.9375
$ws.Range("F21").Formula = "=(E21-D21)*24*60"
$ws.Range("G21").Formula = "=F21/60"

# Move the selection to the newly added time-spent cell.
$ws.Range("F21").Select()
